# Update "想去人数" (interest count) figures in both the "展览" sheet and the
# "全部类型" sheet, which both list the same events and therefore need the
# same updated numbers.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row => new value, for the "展览" sheet (column F)
$sheetExhibition.Range("F3").Value  = 11934
$sheetExhibition.Range("F5").Value  = 222
$sheetExhibition.Range("F6").Value  = 358
$sheetExhibition.Range("F8").Value  = 11836
$sheetExhibition.Range("F9").Value  = 494
$sheetExhibition.Range("F14").Value = 5865

# Same events, same updated values, on the "全部类型" sheet (column F)
$sheetAllTypes.Range("F5").Value   = 11934
$sheetAllTypes.Range("F7").Value   = 222
$sheetAllTypes.Range("F9").Value   = 358
$sheetAllTypes.Range("F11").Value  = 11836
$sheetAllTypes.Range("F12").Value  = 494
$sheetAllTypes.Range("F18").Value  = 5865
